# Updates Faerie Profits figures (columns H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# sheets, per the scheduled-runner commit. Values were recomputed upstream; this script
# just pokes the new numbers into the same cells (adding/clearing cells where the row
# gained or lost a trailing N/M column entry).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2714.4119
$ws.Range("I18").Value = 1305
$ws.Range("K18").Value = 1305
$ws.Range("M18").Value = -1021
$ws.Range("H40").Value = 7147500
$ws.Range("J40").Value = 5166.5835
$ws.Range("L40").Value = 5166.5835
$ws.Range("N40").Value = -5516.5835
$ws.Range("H41").Value = 1212.375
$ws.Range("I41").Value = 939.8
$ws.Range("J41").Value = 1666.6666
$ws.Range("K41").Value = 939.8
$ws.Range("L41").Value = 1666.6666
$ws.Range("M41").Value = -499.8
$ws.Range("N41").Value = -2546.6666
$ws.Range("H53").Value = 1285.1666
$ws.Range("I53").Value = 242.4
$ws.Range("J53").Value = 2030
$ws.Range("K53").Value = 242.4
$ws.Range("L53").Value = 2030
$ws.Range("M53").Value = 394.6
$ws.Range("N53").Value = -3304
$ws.Range("H62").Value = 168031.67
$ws.Range("I62").Value = 201238
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 201238
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -200614
$ws.Range("N62").Value = -3248
$ws.Range("H65").Value = 168031.67
$ws.Range("I65").Value = 201238
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 1006190
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -1003070
$ws.Range("N65").Value = -16240
$ws.Range("H86").Value = 5095.05
$ws.Range("I86").Value = 4422.8184
$ws.Range("K86").Value = 4422.8184
$ws.Range("M86").Value = -3299.8184
$ws.Range("H89").Value = 5095.05
$ws.Range("I89").Value = 4422.8184
$ws.Range("K89").Value = 22114.092
$ws.Range("M89").Value = -16498.092
$ws.Range("H112").Value = 478262.94
$ws.Range("J112").Value = 557651.75
$ws.Range("L112").Value = 1672955.25
$ws.Range("N112").Value = -1675171.25
$ws.Range("H125").Value = 4033.8
$ws.Range("I125").Value = 2204.25
$ws.Range("J125").Value = 5253.5
$ws.Range("K125").Value = 19838.25
$ws.Range("L125").Value = 47281.5
$ws.Range("M125").Value = -17378.25
$ws.Range("N125").Value = -52201.5
$ws.Range("H132").Value = 30309640
$ws.Range("I132").Value = 47626384
$ws.Range("J132").Value = 5334.0835
$ws.Range("K132").Value = 142879152
$ws.Range("L132").Value = 16002.2505
$ws.Range("M132").Value = -142876622
$ws.Range("N132").Value = -21062.2505
$ws.Range("H135").Value = 5659.4165
$ws.Range("I135").Value = 2354
$ws.Range("K135").Value = 21186
$ws.Range("M135").Value = -18651
$ws.Range("H137").Value = 2833.238
$ws.Range("I137").Value = 2189.889
$ws.Range("J137").Value = 3315.75
$ws.Range("K137").Value = 6569.667
$ws.Range("L137").Value = 9947.25
$ws.Range("M137").Value = -4019.667
$ws.Range("N137").Value = -15047.25
$ws.Range("H138").Value = 529093.5
$ws.Range("J138").Value = 1003674.7
$ws.Range("L138").Value = 3011024.1
$ws.Range("N138").Value = -3021304.1

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2238.0527
$ws.Range("I2").Value = 2225.1538
$ws.Range("J2").Value = 2266
$ws.Range("K2").Value = 2225.1538
$ws.Range("L2").Value = 2266
$ws.Range("M2").Value = -2112.1538
$ws.Range("N2").Value = -2492
$ws.Range("H6").Value = 36491
$ws.Range("I6").Value = 60000
$ws.Range("J6").Value = 20818.334
$ws.Range("K6").Value = 60000
$ws.Range("L6").Value = 20818.334
$ws.Range("M6").Value = -59827
$ws.Range("N6").Value = -21164.334
$ws.Range("H25").Value = 175
$ws.Range("I25").Value = 175
$ws.Range("K25").Value = 175
$ws.Range("M25").Value = 227
$ws.Range("H32").Value = 6438.8066
$ws.Range("I32").Value = 8491.736999999999
$ws.Range("J32").Value = 3188.3333
$ws.Range("K32").Value = 8491.736999999999
$ws.Range("L32").Value = 3188.3333
$ws.Range("M32").Value = -8204.736999999999
$ws.Range("N32").Value = -3762.3333
$ws.Range("H37").Value = 23344.666
$ws.Range("I37").Value = 23344.666
$ws.Range("K37").Value = 23344.666
$ws.Range("M37").Value = -23071.666
$ws.Range("H44").Value = 48000
$ws.Range("J44").Value = 48000
$ws.Range("L44").Value = 48000
$ws.Range("N44").Value = -48976
$ws.Range("H49").Value = 40000
$ws.Range("J49").Value = 40000
$ws.Range("L49").Value = 40000
$ws.Range("N49").Value = -40520
$ws.Range("H61").Value = 3792.8823
$ws.Range("I61").Value = 3341.8865
$ws.Range("J61").Value = 6627.7144
$ws.Range("K61").Value = 3341.8865
$ws.Range("L61").Value = 6627.7144
$ws.Range("M61").Value = -3129.8865
$ws.Range("N61").Value = -7051.7144
$ws.Range("H62").Value = 94747
$ws.Range("J62").Value = 94747
$ws.Range("L62").Value = 94747
$ws.Range("N62").Value = -95995
$ws.Range("H65").Value = 94747
$ws.Range("J65").Value = 94747
$ws.Range("L65").Value = 284241
$ws.Range("N65").Value = -290481
$ws.Range("H74").Value = 1826.8572
$ws.Range("I74").Value = 1215.9
$ws.Range("K74").Value = 1215.9
$ws.Range("M74").Value = -341.9000000000001
$ws.Range("H77").Value = 1826.8572
$ws.Range("I77").Value = 1215.9
$ws.Range("K77").Value = 6079.5
$ws.Range("M77").Value = -1711.5
$ws.Range("H93").Value = 68888
$ws.Range("J93").Value = 68888
$ws.Range("L93").Value = 68888
$ws.Range("N93").Value = -73880
$ws.Range("H109").Value = 80377
$ws.Range("J109").Value = 80377
$ws.Range("L109").Value = 80377
$ws.Range("N109").Value = -83151
$ws.Range("H116").Value = 2238.0527
$ws.Range("I116").Value = 2225.1538
$ws.Range("J116").Value = 2266
$ws.Range("K116").Value = 2225.1538
$ws.Range("L116").Value = 2266
$ws.Range("M116").Value = 68.84619999999995
$ws.Range("N116").Value = -6854
$ws.Range("H132").Value = 5702.5884
$ws.Range("I132").Value = 3380.3845
$ws.Range("J132").Value = 13249.75
$ws.Range("K132").Value = 10141.1535
$ws.Range("L132").Value = 39749.25
$ws.Range("M132").Value = -7611.1535
$ws.Range("N132").Value = -44809.25
$ws.Range("H136").Value = 3792.8823
$ws.Range("I136").Value = 3341.8865
$ws.Range("J136").Value = 6627.7144
$ws.Range("K136").Value = 10025.6595
$ws.Range("L136").Value = 19883.1432
$ws.Range("M136").Value = -7475.6595
$ws.Range("N136").Value = -24983.1432

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2238.0527
$ws.Range("I3").Value = 2225.1538
$ws.Range("J3").Value = 2266
$ws.Range("K3").Value = 2225.1538
$ws.Range("L3").Value = 2266
$ws.Range("M3").Value = -2111.1538
$ws.Range("N3").Value = -2494
$ws.Range("H40").Value = 65000
$ws.Range("J40").Value = 65000
$ws.Range("L40").Value = 65000
$ws.Range("N40").Value = -65530
$ws.Range("H62").Value = 79999
$ws.Range("J62").Value = 79999
$ws.Range("L62").Value = 79999
$ws.Range("N62").Value = -81371
$ws.Range("H63").Value = 80000
$ws.Range("J63").Value = 80000
$ws.Range("L63").Value = 80000
$ws.Range("N63").Value = -81372
$ws.Range("H65").Value = 79999
$ws.Range("J65").Value = 79999
$ws.Range("L65").Value = 239997
$ws.Range("N65").Value = -246861
$ws.Range("H66").Value = 80000
$ws.Range("J66").Value = 80000
$ws.Range("L66").Value = 240000
$ws.Range("N66").Value = -246864
$ws.Range("H92").Value = 36411.8
$ws.Range("J92").Value = 36411.8
$ws.Range("L92").Value = 36411.8
$ws.Range("N92").Value = -41403.8
$ws.Range("H93").Value = 72499.5
$ws.Range("J93").Value = 70000
$ws.Range("L93").Value = 70000
$ws.Range("N93").Value = -73744
$ws.Range("H94").Value = 1396.7858
$ws.Range("I94").Value = 1550.6
$ws.Range("J94").Value = 1012.25
$ws.Range("K94").Value = 1550.6
$ws.Range("L94").Value = 1012.25
$ws.Range("M94").Value = -1099.6
$ws.Range("N94").Value = -1914.25
$ws.Range("H96").Value = 66899.60000000001
$ws.Range("I96").Value = 62332.668
$ws.Range("J96").Value = 73750
$ws.Range("K96").Value = 62332.668
$ws.Range("L96").Value = 73750
$ws.Range("M96").Value = -59586.668
$ws.Range("N96").Value = -79242
$ws.Range("H97").Value = 17666.334
$ws.Range("J97").Value = 21499.5
$ws.Range("L97").Value = 21499.5
$ws.Range("N97").Value = -23481.5
$ws.Range("H105").Value = 4442.2
$ws.Range("I105").Value = 6163.6665
$ws.Range("J105").Value = 3704.4285
$ws.Range("K105").Value = 6163.6665
$ws.Range("L105").Value = 3704.4285
$ws.Range("M105").Value = -4416.6665
$ws.Range("N105").Value = -7198.4285
$ws.Range("H134").Value = 1989.2444
$ws.Range("I134").Value = 1635.0513
$ws.Range("K134").Value = 4905.1539
$ws.Range("M134").Value = -2370.1539

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 244.4
$ws.Range("I7").Value = 299.2857
$ws.Range("J7").Value = 116.333336
$ws.Range("K7").Value = 299.2857
$ws.Range("L7").Value = 116.333336
$ws.Range("M7").Value = -186.2857
$ws.Range("N7").Value = -342.333336
$ws.Range("H16").Value = 3449.5
$ws.Range("I16").Value = 3449.5
$ws.Range("K16").Value = 3449.5
$ws.Range("M16").Value = -3162.5
$ws.Range("H31").Value = 2808.879
$ws.Range("I31").Value = 1995.6316
$ws.Range("K31").Value = 1995.6316
$ws.Range("M31").Value = -1700.6316
$ws.Range("H34").Value = 2808.879
$ws.Range("I34").Value = 1995.6316
$ws.Range("K34").Value = 1995.6316
$ws.Range("M34").Value = -1793.6316
$ws.Range("H36").Value = 31798
$ws.Range("I36").Value = 4545
$ws.Range("J36").Value = 49966.668
$ws.Range("K36").Value = 4545
$ws.Range("L36").Value = 49966.668
$ws.Range("M36").Value = -4157
$ws.Range("N36").Value = -50742.668
$ws.Range("H40").Value = 31798
$ws.Range("I40").Value = 4545
$ws.Range("J40").Value = 49966.668
$ws.Range("K40").Value = 4545
$ws.Range("L40").Value = 49966.668
$ws.Range("M40").Value = -4385
$ws.Range("N40").Value = -50286.668
$ws.Range("H58").Value = 2086.0557
$ws.Range("I58").Value = 1178
$ws.Range("J58").Value = 2994.111
$ws.Range("K58").Value = 1178
$ws.Range("L58").Value = 2994.111
$ws.Range("M58").Value = -975
$ws.Range("N58").Value = -3400.111
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = ""
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = ""
$ws.Range("H69").Value = 30110.223
$ws.Range("J69").Value = 35832.168
$ws.Range("L69").Value = 35832.168
$ws.Range("N69").Value = -37330.168
$ws.Range("H70").Value = 22374.875
$ws.Range("J70").Value = 22374.875
$ws.Range("L70").Value = 22374.875
$ws.Range("N70").Value = -23004.875
$ws.Range("H72").Value = 30110.223
$ws.Range("J72").Value = 35832.168
$ws.Range("L72").Value = 107496.504
$ws.Range("N72").Value = -114984.504
$ws.Range("H73").Value = 22374.875
$ws.Range("J73").Value = 22374.875
$ws.Range("L73").Value = 22374.875
$ws.Range("N73").Value = -24558.875
$ws.Range("H97").Value = 98588
$ws.Range("J97").Value = 98588
$ws.Range("L97").Value = 98588
$ws.Range("N97").Value = -100570
$ws.Range("H102").Value = 87212.25
$ws.Range("J102").Value = 87212.25
$ws.Range("L102").Value = 87212.25
$ws.Range("N102").Value = -92080.25
$ws.Range("H113").Value = 3449.5
$ws.Range("I113").Value = 3449.5
$ws.Range("K113").Value = 3449.5
$ws.Range("M113").Value = -1279.5
$ws.Range("H134").Value = 1949.5555
$ws.Range("I134").Value = 1729.4
$ws.Range("K134").Value = 5188.200000000001
$ws.Range("M134").Value = -2653.200000000001
$ws.Range("H136").Value = 2086.0557
$ws.Range("I136").Value = 1178
$ws.Range("J136").Value = 2994.111
$ws.Range("K136").Value = 3534
$ws.Range("L136").Value = 8982.332999999999
$ws.Range("M136").Value = -984
$ws.Range("N136").Value = -14082.333

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1843.5625
$ws.Range("I129").Value = 788.9167
$ws.Range("K129").Value = 2366.7501
$ws.Range("M129").Value = 2633.2499
$ws.Range("H132").Value = 2350.1304
$ws.Range("I132").Value = 1248.5
$ws.Range("K132").Value = 11236.5
$ws.Range("M132").Value = -8706.5
$ws.Range("H134").Value = 1928.55
$ws.Range("I134").Value = 1398.4736
$ws.Range("K134").Value = 4195.4208
$ws.Range("M134").Value = 874.5792000000001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 1101.5385
$ws.Range("I13").Value = 308
$ws.Range("K13").Value = 308
$ws.Range("M13").Value = -169
$ws.Range("H62").Value = 79994
$ws.Range("J62").Value = 79994
$ws.Range("L62").Value = 79994
$ws.Range("N62").Value = -81366
$ws.Range("H63").Value = 94999
$ws.Range("J63").Value = 94999
$ws.Range("L63").Value = 94999
$ws.Range("N63").Value = -96371
$ws.Range("H65").Value = 79994
$ws.Range("J65").Value = 79994
$ws.Range("L65").Value = 239982
$ws.Range("N65").Value = -246846
$ws.Range("H66").Value = 94999
$ws.Range("J66").Value = 94999
$ws.Range("L66").Value = 284997
$ws.Range("N66").Value = -291861
$ws.Range("H68").Value = 42498.5
$ws.Range("J68").Value = 42498.5
$ws.Range("L68").Value = 42498.5
$ws.Range("N68").Value = -44120.5
$ws.Range("H69").Value = 69404.664
$ws.Range("J69").Value = 69404.664
$ws.Range("L69").Value = 69404.664
$ws.Range("N69").Value = -70902.664
$ws.Range("H70").Value = 10165.111
$ws.Range("I70").Value = 3747.25
$ws.Range("J70").Value = 15299.4
$ws.Range("K70").Value = 3747.25
$ws.Range("L70").Value = 15299.4
$ws.Range("M70").Value = -3477.25
$ws.Range("N70").Value = -15839.4
$ws.Range("H71").Value = 42498.5
$ws.Range("J71").Value = 42498.5
$ws.Range("L71").Value = 127495.5
$ws.Range("N71").Value = -135607.5
$ws.Range("H72").Value = 69404.664
$ws.Range("J72").Value = 69404.664
$ws.Range("L72").Value = 208213.992
$ws.Range("N72").Value = -215701.992
$ws.Range("H73").Value = 10165.111
$ws.Range("I73").Value = 3747.25
$ws.Range("J73").Value = 15299.4
$ws.Range("K73").Value = 3747.25
$ws.Range("L73").Value = 15299.4
$ws.Range("M73").Value = -2811.25
$ws.Range("N73").Value = -17171.4
$ws.Range("H97").Value = 1762.409
$ws.Range("I97").Value = 994.40625
$ws.Range("K97").Value = 994.40625
$ws.Range("M97").Value = -498.40625
$ws.Range("H99").Value = 30632.75
$ws.Range("I99").Value = 1735.5
$ws.Range("K99").Value = 1735.5
$ws.Range("M99").Value = 510.5
$ws.Range("H126").Value = 19325.5
$ws.Range("I126").Value = 4766.769
$ws.Range("J126").Value = 57178.2
$ws.Range("K126").Value = 14300.307
$ws.Range("L126").Value = 171534.6
$ws.Range("M126").Value = -11830.307
$ws.Range("N126").Value = -176474.6
$ws.Range("H132").Value = 10045.346
$ws.Range("I132").Value = 10551.261
$ws.Range("J132").Value = 6166.6665
$ws.Range("K132").Value = 31653.783
$ws.Range("L132").Value = 18499.9995
$ws.Range("M132").Value = -29123.783
$ws.Range("N132").Value = -23559.9995

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1791.875
$ws.Range("I16").Value = 1226.375
$ws.Range("J16").Value = 2357.375
$ws.Range("K16").Value = 1226.375
$ws.Range("L16").Value = 2357.375
$ws.Range("M16").Value = -1056.375
$ws.Range("N16").Value = -2697.375
$ws.Range("H45").Value = 29041
$ws.Range("I45").Value = 29041
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 29041
$ws.Range("L45").Value = 0
$ws.Range("N45").Value = ""
$ws.Range("M45").Value = -28634
$ws.Range("H46").Value = 27700.445
$ws.Range("J46").Value = 43661.6
$ws.Range("L46").Value = 43661.6
$ws.Range("N46").Value = -44037.6
$ws.Range("H96").Value = 89898
$ws.Range("J96").Value = 89898
$ws.Range("L96").Value = 89898
$ws.Range("N96").Value = -95390
$ws.Range("H99").Value = 62086.332
$ws.Range("I99").Value = 21259
$ws.Range("J99").Value = 82500
$ws.Range("K99").Value = 21259
$ws.Range("L99").Value = 82500
$ws.Range("N99").Value = -88490
$ws.Range("M99").Value = -18264
$ws.Range("H100").Value = 3590.5833
$ws.Range("I100").Value = 3376.4443
$ws.Range("J100").Value = 4233
$ws.Range("K100").Value = 3376.4443
$ws.Range("L100").Value = 4233
$ws.Range("M100").Value = -2835.4443
$ws.Range("N100").Value = -5315
$ws.Range("H102").Value = 96779.5
$ws.Range("J102").Value = 96779.5
$ws.Range("L102").Value = 96779.5
$ws.Range("N102").Value = -103269.5
$ws.Range("H132").Value = 5827.696
$ws.Range("I132").Value = 5126.3076
$ws.Range("K132").Value = 15378.9228
$ws.Range("M132").Value = -12848.9228
$ws.Range("H136").Value = 3507.2444
$ws.Range("I136").Value = 3426.3333
$ws.Range("J136").Value = 4033.1667
$ws.Range("K136").Value = 10278.9999
$ws.Range("L136").Value = 12099.5001
$ws.Range("M136").Value = -7728.999899999999
$ws.Range("N136").Value = -17199.5001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 1682333.4
$ws.Range("I5").Value = 3339666.8
$ws.Range("J5").Value = 25000
$ws.Range("K5").Value = 3339666.8
$ws.Range("L5").Value = 25000
$ws.Range("M5").Value = -3339554.8
$ws.Range("N5").Value = -25224
$ws.Range("H6").Value = 1001.875
$ws.Range("I6").Value = 176.66667
$ws.Range("J6").Value = 1497
$ws.Range("K6").Value = 176.66667
$ws.Range("L6").Value = 1497
$ws.Range("M6").Value = -61.66667000000001
$ws.Range("N6").Value = -1727
$ws.Range("H11").Value = 2510000
$ws.Range("I11").Value = 2510000
$ws.Range("K11").Value = 2510000
$ws.Range("M11").Value = -2509858
$ws.Range("H13").Value = 5499.5
$ws.Range("I13").Value = 5499.5
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 5499.5
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -5359.5
$ws.Range("N13").Value = ""
$ws.Range("H17").Value = 5600
$ws.Range("I17").Value = 5600
$ws.Range("K17").Value = 5600
$ws.Range("M17").Value = -5428
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = ""
$ws.Range("N18").Value = ""
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").Value = ""
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = ""
$ws.Range("H23").Value = 1014.2222
$ws.Range("I23").Value = 1176.8572
$ws.Range("J23").Value = 445
$ws.Range("K23").Value = 1176.8572
$ws.Range("L23").Value = 445
$ws.Range("M23").Value = -947.8571999999999
$ws.Range("N23").Value = -903
$ws.Range("H24").Value = 10000
$ws.Range("I24").Value = 10000
$ws.Range("K24").Value = 10000
$ws.Range("M24").Value = -9770
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").Value = ""
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").Value = ""
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").Value = ""
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").Value = ""
$ws.Range("H39").Value = 2509500
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 2509500
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 2509500
$ws.Range("M39").Value = ""
$ws.Range("N39").Value = -2510326
$ws.Range("H41").Value = 22263.143
$ws.Range("I41").Value = 23671
$ws.Range("J41").Value = 21700
$ws.Range("K41").Value = 23671
$ws.Range("L41").Value = 21700
$ws.Range("M41").Value = -23281
$ws.Range("N41").Value = -22480
$ws.Range("H44").Value = 30000
$ws.Range("J44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -31108
$ws.Range("H48").Value = 31583
$ws.Range("I48").Value = 29750
$ws.Range("J48").Value = 32499.5
$ws.Range("K48").Value = 29750
$ws.Range("L48").Value = 32499.5
$ws.Range("M48").Value = -29181
$ws.Range("N48").Value = -33637.5
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").Value = ""
$ws.Range("H54").Value = 20500
$ws.Range("J54").Value = 20500
$ws.Range("L54").Value = 20500
$ws.Range("N54").Value = -21540
$ws.Range("H59").Value = 42221
$ws.Range("J59").Value = 42221
$ws.Range("L59").Value = 42221
$ws.Range("N59").Value = -43697
$ws.Range("H61").Value = 19666.334
$ws.Range("J61").Value = 19999.5
$ws.Range("L61").Value = 19999.5
$ws.Range("N61").Value = -20583.5
$ws.Range("H62").Value = 263110.66
$ws.Range("I62").Value = 391888.5
$ws.Range("K62").Value = 391888.5
$ws.Range("M62").Value = -391264.5
$ws.Range("H63").Value = 47500
$ws.Range("J63").Value = 47500
$ws.Range("L63").Value = 47500
$ws.Range("N63").Value = -48748
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = ""
$ws.Range("H65").Value = 263110.66
$ws.Range("I65").Value = 391888.5
$ws.Range("K65").Value = 1959442.5
$ws.Range("M65").Value = -1956322.5
$ws.Range("H66").Value = 47500
$ws.Range("J66").Value = 47500
$ws.Range("L66").Value = 142500
$ws.Range("N66").Value = -148740
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = ""
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = ""
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = ""
$ws.Range("H87").Value = 96420
$ws.Range("J87").Value = 96420
$ws.Range("L87").Value = 96420
$ws.Range("N87").Value = -98916
$ws.Range("H90").Value = 96420
$ws.Range("J90").Value = 96420
$ws.Range("L90").Value = 289260
$ws.Range("N90").Value = -301740
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""
$ws.Range("H93").Value = 88604.5
$ws.Range("J93").Value = 88604.5
$ws.Range("L93").Value = 88604.5
$ws.Range("N93").Value = -93596.5
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = ""
$ws.Range("H96").Value = 4071.375
$ws.Range("I96").Value = 3643
$ws.Range("K96").Value = 3643
$ws.Range("M96").Value = -2270
$ws.Range("H97").Value = 69999
$ws.Range("J97").Value = 69999
$ws.Range("L97").Value = 69999
$ws.Range("N97").Value = -71981
$ws.Range("H98").Value = 80000
$ws.Range("J98").Value = 80000
$ws.Range("L98").Value = 80000
$ws.Range("N98").Value = -85990
$ws.Range("H102").Value = 92949
$ws.Range("J102").Value = 92949
$ws.Range("L102").Value = 92949
$ws.Range("N102").Value = -99439
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = ""
$ws.Range("H106").Value = 87688.5
$ws.Range("J106").Value = 87688.5
$ws.Range("L106").Value = 87688.5
$ws.Range("N106").Value = -90212.5
$ws.Range("H116").Value = 85681
$ws.Range("J116").Value = 85681
$ws.Range("L116").Value = 85681
$ws.Range("N116").Value = -94859
$ws.Range("H120").Value = 91909.5
$ws.Range("J120").Value = 91909.5
$ws.Range("L120").Value = 91909.5
$ws.Range("N120").Value = -101585.5
$ws.Range("H122").Value = 6846.919
$ws.Range("I122").Value = 7060.4614
$ws.Range("K122").Value = 21181.3842
$ws.Range("M122").Value = -18731.3842
$ws.Range("H126").Value = 8582.322
$ws.Range("I126").Value = 7980.6523
$ws.Range("K126").Value = 23941.9569
$ws.Range("M126").Value = -21471.9569
$ws.Range("H136").Value = 2862.5789
$ws.Range("I136").Value = 2687.4375
$ws.Range("K136").Value = 8062.3125
$ws.Range("M136").Value = -5512.3125

Write-Host "Applied 660 cell updates"
